$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.050.33"
$ws.Range("E2").Value = "  +5.87%  "

$ws.Range("D3").Value = "2.429.69"
$ws.Range("E3").Value = "  +5.37%  "

$ws.Range("D5").Value = "'564.22"
$ws.Range("E5").Value = "  +4.35%  "

$ws.Range("D6").Value = "'141.84"
$ws.Range("E6").Value = "  +11.26%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  +3.33%  "

$ws.Range("D9").Value = "2.429.77"
$ws.Range("E9").Value = "  +5.51%  "

$ws.Range("E10").Value = "  +4.21%  "

$ws.Range("D11").Value = "'5.74"
$ws.Range("E11").Value = "  +4.09%  "

$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  +6.32%  "

$ws.Range("D14").Value = "'26.33"
$ws.Range("E14").Value = "  +14.21%  "

$ws.Range("D15").Value = "2.866.17"
$ws.Range("E15").Value = "  +5.51%  "

$ws.Range("D16").Value = "62.937.53"
$ws.Range("E16").Value = "  +5.67%  "

$ws.Range("E17").Value = "  +8.66%  "

$ws.Range("D18").Value = "2.429.81"
$ws.Range("E18").Value = "  +5.12%  "

$ws.Range("D19").Value = "'11.20"
$ws.Range("E19").Value = "  +7.73%  "

$ws.Range("D20").Value = "'338.97"
$ws.Range("E20").Value = "  +9.47%  "

$ws.Range("E21").Value = "  +5.93%  "

$ws.Range("D22").Value = "'6.76"
$ws.Range("E22").Value = "  +3.99%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'65.26"
$ws.Range("E24").Value = "  +3.61%  "

$ws.Range("D25").Value = "'0.173"
$ws.Range("E25").Value = "  +3.50%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  +14.05%  "

$ws.Range("D28").Value = "'8.13"
$ws.Range("E28").Value = "  +5.82%  "

$ws.Range("D29").Value = "'1.32"
$ws.Range("E29").Value = "  +11.83%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'6.67"
$ws.Range("E30").Value = "  +15.51%  "

$ws.Range("E31").Value = "  +7.22%  "

$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0790"
$ws.Range("E32").Value = "  +10.80%  "

$ws.Range("D33").Value = "'174.77"
$ws.Range("E33").Value = "  +1.78%  "

$ws.Range("E34").Value = "  +11.93%  "

$ws.Range("E35").Value = "  +6.12%  "

$ws.Range("E36").Value = "  +5.84%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = "  +12.87%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'368.59"
$ws.Range("E38").Value = "  +18.21%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").Value = "'1.71"
$ws.Range("E41").Value = "  +14.27%  "

$ws.Range("D42").Value = "'40.07"
$ws.Range("E42").Value = "  +6.69%  "

$ws.Range("D43").Value = "'149.55"
$ws.Range("E43").Value = "  +9.96%  "

$ws.Range("E44").Value = "  +8.50%  "

$ws.Range("D45").Value = "'20.65"
$ws.Range("E45").Value = "  +12.09%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.593"
$ws.Range("E46").Value = "  +4.45%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0957"
$ws.Range("E47").Value = "  +2.25%  "

$ws.Range("D48").Value = "'0.0521"
$ws.Range("E48").Value = "  +6.79%  "

$ws.Range("D49").Value = "'0.0225"
$ws.Range("E49").Value = "  +6.71%  "

$ws.Range("D50").Value = "'17.80"
$ws.Range("E50").Value = "  +6.94%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0223"
$ws.Range("E51").Value = "  +0.09%  "
